$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws1.Range("H28").Value = 1122.8889
$ws1.Range("I28").Value = 1122.8889
$ws1.Range("K28").Value = 1122.8889
$ws1.Range("M28").Value = -637.8888999999999
$ws1.Range("H86").Value = 4854.643
$ws1.Range("I86").Value = 3399.6667
$ws1.Range("J86").Value = 5251.4546
$ws1.Range("K86").Value = 3399.6667
$ws1.Range("L86").Value = 5251.4546
$ws1.Range("M86").Value = -2276.6667
$ws1.Range("N86").Value = -7497.4546
$ws1.Range("H89").Value = 4854.643
$ws1.Range("I89").Value = 3399.6667
$ws1.Range("J89").Value = 5251.4546
$ws1.Range("K89").Value = 16998.3335
$ws1.Range("L89").Value = 26257.273
$ws1.Range("M89").Value = -11382.3335
$ws1.Range("N89").Value = -37489.273
$ws1.Range("H92").Value = 643.8
$ws1.Range("I92").Value = 579.875
$ws1.Range("K92").Value = 579.875
$ws1.Range("M92").Value = 668.125
$ws1.Range("H94").Value = 673.5
$ws1.Range("I94").Value = 673.5
$ws1.Range("K94").Value = 673.5
$ws1.Range("M94").Value = -222.5
$ws1.Range("H103").Value = 3004.75
$ws1.Range("I103").Value = 2763.25
$ws1.Range("J103").Value = 3246.25
$ws1.Range("K103").Value = 8289.75
$ws1.Range("L103").Value = 9738.75
$ws1.Range("M103").Value = -7703.75
$ws1.Range("N103").Value = -10910.75
$ws1.Range("H107").Value = 545.6
$ws1.Range("I107").Value = 676.6667
$ws1.Range("K107").Value = 676.6667
$ws1.Range("M107").Value = 1243.3333
$ws1.Range("H132").Value = 965.3043
$ws1.Range("I132").Value = 990.15
$ws1.Range("K132").Value = 2970.45
$ws1.Range("M132").Value = -440.4499999999998
$ws1.Range("H135").Value = 1111.3077
$ws1.Range("I135").Value = 1129.9166
$ws1.Range("K135").Value = 10169.2494
$ws1.Range("M135").Value = -7634.249400000001
$ws1.Range("H137").Value = 1216.4445
$ws1.Range("I137").Value = 931.125
$ws1.Range("K137").Value = 2793.375
$ws1.Range("M137").Value = -243.375
$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws2.Range("H2").Value = 1774.5
$ws2.Range("I2").Value = 762
$ws2.Range("J2").Value = 3799.5
$ws2.Range("K2").Value = 762
$ws2.Range("L2").Value = 3799.5
$ws2.Range("M2").Value = -649
$ws2.Range("N2").Value = -4025.5
$ws2.Range("H32").Value = 5297.968
$ws2.Range("I32").Value = 5008.1724
$ws2.Range("K32").Value = 5008.1724
$ws2.Range("M32").Value = -4721.1724
$ws2.Range("H74").Value = 2724.353
$ws2.Range("I74").Value = 2624.0833
$ws2.Range("J74").Value = 2965
$ws2.Range("K74").Value = 2624.0833
$ws2.Range("L74").Value = 2965
$ws2.Range("M74").Value = -1750.0833
$ws2.Range("N74").Value = -4713
$ws2.Range("H77").Value = 2724.353
$ws2.Range("I77").Value = 2624.0833
$ws2.Range("J77").Value = 2965
$ws2.Range("K77").Value = 13120.4165
$ws2.Range("L77").Value = 14825
$ws2.Range("M77").Value = -8752.416499999999
$ws2.Range("N77").Value = -23561
$ws2.Range("H97").Value = 822.5714
$ws2.Range("I97").Value = 843
$ws2.Range("K97").Value = 843
$ws2.Range("M97").Value = -347
$ws2.Range("H116").Value = 1774.5
$ws2.Range("I116").Value = 762
$ws2.Range("J116").Value = 3799.5
$ws2.Range("K116").Value = 762
$ws2.Range("L116").Value = 3799.5
$ws2.Range("M116").Value = 1532
$ws2.Range("N116").Value = -8387.5
$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws3.Range("H3").Value = 1774.5
$ws3.Range("I3").Value = 762
$ws3.Range("J3").Value = 3799.5
$ws3.Range("K3").Value = 762
$ws3.Range("L3").Value = 3799.5
$ws3.Range("M3").Value = -648
$ws3.Range("N3").Value = -4027.5
$ws3.Range("H107").Value = 2672.3157
$ws3.Range("I107").Value = 2375.611
$ws3.Range("J107").Value = 8013
$ws3.Range("K107").Value = 2375.611
$ws3.Range("L107").Value = 8013
$ws3.Range("M107").Value = -455.6109999999999
$ws3.Range("N107").Value = -11853
$ws3.Range("H134").Value = 488
$ws3.Range("I134").Value = 488
$ws3.Range("K134").Value = 1464
$ws3.Range("M134").Value = 1071
$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws4.Range("H16").Value = 2466.625
$ws4.Range("I16").Value = 2117.1428
$ws4.Range("J16").Value = 4913
$ws4.Range("K16").Value = 2117.1428
$ws4.Range("L16").Value = 4913
$ws4.Range("M16").Value = -1830.1428
$ws4.Range("N16").Value = -5487
$ws4.Range("H107").Value = 619.7692
$ws4.Range("I107").Value = 497.33334
$ws4.Range("J107").Value = 895.25
$ws4.Range("K107").Value = 497.33334
$ws4.Range("L107").Value = 895.25
$ws4.Range("M107").Value = 1422.66666
$ws4.Range("N107").Value = -4735.25
$ws4.Range("H113").Value = 2466.625
$ws4.Range("I113").Value = 2117.1428
$ws4.Range("J113").Value = 4913
$ws4.Range("K113").Value = 2117.1428
$ws4.Range("L113").Value = 4913
$ws4.Range("M113").Value = 52.85719999999992
$ws4.Range("N113").Value = -9253
$ws4.Range("H134").Value = 1708.0667
$ws4.Range("I134").Value = 1708.0667
$ws4.Range("K134").Value = 5124.2001
$ws4.Range("M134").Value = -2589.2001
$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws5.Range("H4").Value = 15714457
$ws5.Range("I4").Value = 15714457
$ws5.Range("K4").Value = 47143371
$ws5.Range("M4").Value = -47143259
$ws5.Range("H56").Value = 8722.24
$ws5.Range("I56").Value = 8722.24
$ws5.Range("K56").Value = 8722.24
$ws5.Range("M56").Value = -8192.24
$ws5.Range("H98").Value = 238.2
$ws5.Range("I98").Value = 224.25
$ws5.Range("J98").Value = 294
$ws5.Range("K98").Value = 672.75
$ws5.Range("L98").Value = 882
$ws5.Range("M98").Value = 825.25
$ws5.Range("N98").Value = -3878
$ws5.Range("H140").Value = 10647
$ws5.Range("I140").Value = 7382.25
$ws5.Range("K140").Value = 22146.75
$ws5.Range("M140").Value = -16966.75
$ws7 = $wb.Worksheets.Item(7)  # LTW
$ws7.Range("H46").Value = 2769.3572
$ws7.Range("I46").Value = 1267.5
$ws7.Range("K46").Value = 1267.5
$ws7.Range("M46").Value = -1079.5
$ws7.Range("H122").Value = 6586.0884
$ws7.Range("I122").Value = 6437.68
$ws7.Range("K122").Value = 19313.04
$ws7.Range("M122").Value = -16863.04
$ws8 = $wb.Worksheets.Item(8)  # WVR
$ws8.Range("H107").Value = 577.5
$ws8.Range("I107").Value = 555
$ws8.Range("J107").Value = 600
$ws8.Range("K107").Value = 1665
$ws8.Range("L107").Value = 1800
$ws8.Range("M107").Value = 255
$ws8.Range("N107").Value = -5640
$ws8.Range("H136").Value = 3545.889
$ws8.Range("I136").Value = 3530.8
$ws8.Range("K136").Value = 10592.4
$ws8.Range("M136").Value = -8042.400000000001
